# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (column G) of the attendance report.
#
# Original values like "System, dnasr281@gmail.com" become
# "dnasr281@gmail.com, System". Cells that already contain only one of
# the two tokens (e.g. just "System" or just an email) are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
